$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("NewLoanInput")
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# --- Repayment Schedule sheet: clear stale zero placeholders & drop the
# empty "O" separator column before the tab/selection change below so the
# Repayment Schedule is left in the state it was in just before the user
# tabbed back over to NewLoanInput.

# Row 2 holds several explicit 0 placeholders that become blank once the
# loan figures are recalculated for the new (multi-browser) run.
$row2Cols = @("A", "B", "D", "E", "F", "H", "J", "M", "N", "P")
foreach ($col in $row2Cols) {
    $wsSchedule.Range($col + "2").ClearContents()
}

# Rows 3-14: the "Paid Date" / blank heading columns (D, E) clear out too.
for ($r = 3; $r -le 14; $r++) {
    $wsSchedule.Range("D" + $r).ClearContents()
    $wsSchedule.Range("E" + $r).ClearContents()
}

# Column O (an unlabeled "heading" separator column) is removed outright
# for every data row, not just blanked.
for ($r = 2; $r -le 14; $r++) {
    $wsSchedule.Range("O" + $r).Clear()
}

# Move the Repayment Schedule selection, then leave it as the inactive tab.
$wsSchedule.Range("F20").Select() | Out-Null

# NewLoanInput becomes the active/selected tab with a new selection.
$wsInput.Select() | Out-Null
$wsInput.Range("B22").Select() | Out-Null
